$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EngimonList")

$icons = @(
    "Icons/icon053.png",
    "Icons/icon007.png",
    "Icons/icon077.png",
    "Icons/icon148.png",
    "Icons/icon350.png",
    "Icons/icon446.png",
    "Icons/icon184.png",
    "Icons/icon392.png",
    "Icons/icon485.png",
    "Icons/icon479.png",
    "Icons/icon373.png",
    "Icons/icon376.png",
    "Icons/icon378.png",
    "Icons/icon381.png",
    "Icons/icon386.png",
    "Icons/icon310.png",
    "Icons/icon145.png",
    "Icons/icon160.png",
    "Icons/icon257.png",
    "Icons/icon275.png",
    "Icons/icon375.png",
    "Icons/icon377.png",
    "Icons/icon384.png",
    "Icons/icon488.png",
    "Icons/icon486.png",
    "Icons/icon623.png",
    "Icons/icon382.png",
    "Icons/icon643.png",
    "Icons/icon448.png",
    "Icons/icon445.png",
    "Icons/icon000.png",
    "Icons/icon244.png",
    "Icons/icon065.png"
)

for ($i = 0; $i -lt $icons.Length; $i++) {
    $row = $i + 5
    $ws.Cells.Item($row, 7).Value = $icons[$i]
}

$ws.Range("G38").Select()
